$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 68
$ws.Range("I2").Value = 185
$ws.Range("J2").Value = 695
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 186
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = 133
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 11
$ws.Range("S2").Value = 80
$ws.Range("T2").Value = 122
$ws.Range("U2").Value = 8
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1111
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 16
$ws.Range("AA2").Value = 8
